$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    ,@(42,160,150,3,1,13,3,44,1,0,325000)
    ,@(43,170,140,3,1,18,4,4,1,1,277000)
    ,@(44,350,260,5,1,18,3,4,1,0,410000)
    ,@(45,135,105,2,1,28,1,3,1,0,227000)
    ,@(46,145,125,3,1,1,2,4,1,0,314000)
    ,@(47,180,150,3,1,28,5,5,2,0,325000)
    ,@(48,160,140,3,1,28,5,5,1,0,260000)
    ,@(49,115,95,2,1,0,2,4,1,0,280000)
    ,@(50,110,100,2,1,0,2,3,1,0,310000)
    ,@(51,400,350,6,1,13,4,5,1,0,435000)
    ,@(52,120,95,2,1,0,0,4,1,0,199000)
    ,@(53,150,130,3,1,28,3,5,1,0,235000)
    ,@(54,120,100,2,1,35,5,6,1,0,140000)
    ,@(55,170,145,3,1,13,0,5,1,0,327000)
    ,@(56,160,145,3,1,0,2,4,1,0,480000)
    ,@(57,147,137,3,1,8,0,13,1,1,170000)
    ,@(58,164,140,3,1,28,2,5,3,0,169900)
    ,@(59,115,100,2,1,23,1,4,3,0,150000)
    ,@(60,150,125,3,1,35,4,5,3,0,137000)
    ,@(61,120,112,2,1,28,4,6,3,0,138000)
)

foreach ($row in $newData) {
    $r = $row[0]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $row[$c]
    }
    $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,9)).Style = "Giriş"
    $ws.Range($ws.Cells.Item($r,10), $ws.Cells.Item($r,10)).Style = "Çıkış"
}

$ws.Rows.Item(42).RowHeight = 14.25
$ws.Rows.Item(43).RowHeight = 14.25

[void]$ws.Range("N50").Select()
